# Applies the weekly fruit/hortaliza price update (Chirimoya, Agro Chillan).
# Rows 2,3,6,7,8,9,10,11,12,13,14,15,16,17 get their
# Fecha (D), Calidad (L), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) columns reassigned
# to a new set of values (rows 4 and 5 are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $D, $L, $M, $N, $O, $P, $S) {
    $ws.Cells.Item($Row, 4).Value2 = $D    # D: Fecha
    $ws.Cells.Item($Row, 12).Value = $L    # L: Calidad
    $ws.Cells.Item($Row, 13).Value2 = $M   # M: Volumen
    $ws.Cells.Item($Row, 14).Value2 = $N   # N: Precio minimo
    $ws.Cells.Item($Row, 15).Value2 = $O   # O: Precio maximo
    $ws.Cells.Item($Row, 16).Value2 = $P   # P: Precio promedio ponderado
    $ws.Cells.Item($Row, 19).Value2 = $S   # S: Precio $/Kg
}

Set-Row 2  44839 "Primera"  120 25000 26000 25500 2550
Set-Row 3  44841 "Primera"  60  23000 24000 23500 2350
Set-Row 6  44848 "Especial" 60  24000 25000 24500 2450
Set-Row 7  44848 "Primera"  120 21000 22000 21500 2150
Set-Row 8  44487 "Primera"  30  23000 24000 23500 2350
Set-Row 9  44452 "Primera"  60  21000 22000 21500 2150
Set-Row 10 44448 "Primera"  60  21000 22000 21500 2150
Set-Row 11 44461 "Especial" 60  31000 32000 31500 3150
Set-Row 12 44461 "Primera"  30  30000 30000 30000 3000
Set-Row 13 44447 "Primera"  60  21000 22000 21500 2150
Set-Row 14 44868 "Especial" 60  26000 26000 26000 2600
Set-Row 15 44874 "Especial" 30  25000 25000 25000 2500
Set-Row 16 44874 "Primera"  80  23000 24000 23500 2350
Set-Row 17 44446 "Primera"  60  21000 22000 21500 2150
